# Add a new "JoinCastleItemID" (回城卷轴ID) column (K) to the ConstAttConfig
# sheet: header / type / description rows plus the data value 112.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ConstAttConfig")

# K3 (description row) uses the same "bold-ish" style as the other description
# cells in row 3 (e.g. B3/D3/G3/H3/J3) - copy that formatting over first.
$ws.Range("B3").Copy()
$ws.Range("K3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("K1").Value = "JoinCastleItemID"
$ws.Range("K2").Value = "Int"
$ws.Range("K3").Value = "回城卷轴ID"
$ws.Range("K4").Value = 112

# Widen the new column to match the source width (~18.625 characters).
$ws.Columns.Item(11).ColumnWidth = (125/7)

# Leave the selection on the new column, as in the authored edit.
$ws.Range("H13").Select()
